# "second max element and started the searching in arrays by linearch
#  searchin arrays" -- the author is making room at the bottom of the
# Array Problems document for the next couple of problems (second max
# element, linear search) that are about to be typed in. Mechanically,
# that shows up here as two new blank paragraphs appended right at the
# end of the document, just ahead of the already-present trailing blank
# paragraph / section break.

$d = $word.ActiveDocument

# Collapse a range to the very end of the document content, then insert
# two blank paragraph marks before that point (mirrors a user placing the
# cursor at the end of the doc and pressing Enter twice).
$endRange = $d.Content
$endRange.Collapse(0)   # wdCollapseEnd
$endRange.InsertParagraphBefore()
$endRange.InsertParagraphBefore()
